# Auto-generated edit script: updates market-price-derived columns (H-N)
# on several leve-profit sheets, refreshed by the scheduled pricing runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2897.5
$ws.Range("I17").Value = 2833.3333
$ws.Range("J17").Value = 3090
$ws.Range("K17").Value = 8499.999899999999
$ws.Range("L17").Value = 9270
$ws.Range("M17").Value = -8331.999899999999
$ws.Range("N17").Value = -9606
$ws.Range("H18").Value = 800
$ws.Range("I18").Value = 800
$ws.Range("J18").Value = 800
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = -516
$ws.Range("N18").Value = -1368
$ws.Range("H70").Value = 6143.0835
$ws.Range("J70").Value = 9843.4
$ws.Range("L70").Value = 29530.2
$ws.Range("N70").Value = -30070.2
$ws.Range("H73").Value = 6143.0835
$ws.Range("J73").Value = 9843.4
$ws.Range("L73").Value = 29530.2
$ws.Range("N73").Value = -31402.2
$ws.Range("H80").Value = 2258
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 15000
$ws.Range("M80").Value = -14002
$ws.Range("H83").Value = 2258
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40008
$ws.Range("H98").Value = 19399.3
$ws.Range("I98").Value = 17110.334
$ws.Range("K98").Value = 17110.334
$ws.Range("M98").Value = -15612.334
$ws.Range("H99").Value = 916.4545000000001
$ws.Range("I99").Value = 291.66666
$ws.Range("J99").Value = 1150.75
$ws.Range("K99").Value = 874.9999799999999
$ws.Range("L99").Value = 3452.25
$ws.Range("M99").Value = 623.0000200000001
$ws.Range("N99").Value = -6448.25
$ws.Range("H106").Value = 707.2
$ws.Range("I106").Value = 707.2
$ws.Range("K106").Value = 707.2
$ws.Range("M106").Value = -76.20000000000005
$ws.Range("H112").Value = 5824.75
$ws.Range("I112").Value = 3250
$ws.Range("J112").Value = 8399.5
$ws.Range("K112").Value = 9750
$ws.Range("L112").Value = 25198.5
$ws.Range("M112").Value = -8642
$ws.Range("N112").Value = -27414.5
$ws.Range("H113").Value = 11369.818
$ws.Range("I113").Value = 10025
$ws.Range("J113").Value = 12983.6
$ws.Range("K113").Value = 10025
$ws.Range("L113").Value = 12983.6
$ws.Range("M113").Value = -6771
$ws.Range("N113").Value = -19491.6
$ws.Range("H122").Value = 19399.3
$ws.Range("I122").Value = 17110.334
$ws.Range("K122").Value = 51331.00199999999
$ws.Range("M122").Value = -48881.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1199
$ws.Range("I61").Value = 1199
$ws.Range("K61").Value = 1199
$ws.Range("M61").Value = -987
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344
$ws.Range("H136").Value = 1199
$ws.Range("I136").Value = 1199
$ws.Range("K136").Value = 3597
$ws.Range("M136").Value = -1047

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 50010.5
$ws.Range("I33").Value = 20021
$ws.Range("K33").Value = 20021
$ws.Range("M33").Value = -19685

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 449.5
$ws.Range("I16").Value = 449.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 449.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -162.5
$ws.Range("N16").Value = ""
$ws.Range("H22").Value = 277
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H99").Value = 1251749.8
$ws.Range("I99").Value = 2499.5
$ws.Range("J99").Value = 2501000
$ws.Range("K99").Value = 2499.5
$ws.Range("L99").Value = 2501000
$ws.Range("M99").Value = -1001.5
$ws.Range("N99").Value = -2503996
$ws.Range("H107").Value = 617.9091
$ws.Range("I107").Value = 582.5
$ws.Range("J107").Value = 660.4
$ws.Range("K107").Value = 582.5
$ws.Range("L107").Value = 660.4
$ws.Range("M107").Value = 1337.5
$ws.Range("N107").Value = -4500.4
$ws.Range("H113").Value = 449.5
$ws.Range("I113").Value = 449.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 449.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1720.5
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 3075.4285
$ws.Range("I122").Value = 213.5
$ws.Range("K122").Value = 640.5
$ws.Range("M122").Value = 1809.5
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -79920
$ws.Range("H126").Value = 1251749.8
$ws.Range("I126").Value = 2499.5
$ws.Range("J126").Value = 2501000
$ws.Range("K126").Value = 7498.5
$ws.Range("L126").Value = 7503000
$ws.Range("M126").Value = -5028.5
$ws.Range("N126").Value = -7507940
$ws.Range("H132").Value = 2018.5
$ws.Range("I132").Value = 1909.4445
$ws.Range("K132").Value = 5728.333500000001
$ws.Range("M132").Value = -3198.333500000001
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = ""
$ws.Range("H113").Value = 1699.3334
$ws.Range("I113").Value = 1699.3334
$ws.Range("K113").Value = 1699.3334
$ws.Range("M113").Value = 470.6666
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 704.5
$ws.Range("I9").Value = 410
$ws.Range("J9").Value = 999
$ws.Range("K9").Value = 410
$ws.Range("L9").Value = 999
$ws.Range("M9").Value = -186
$ws.Range("N9").Value = -1447
$ws.Range("H93").Value = 1466.1666
$ws.Range("I93").Value = 1416
$ws.Range("J93").Value = 1516.3334
$ws.Range("K93").Value = 1416
$ws.Range("L93").Value = 1516.3334
$ws.Range("M93").Value = -168
$ws.Range("N93").Value = -4012.3334
$ws.Range("H100").Value = 3402.5715
$ws.Range("I100").Value = 2863.6
$ws.Range("K100").Value = 2863.6
$ws.Range("M100").Value = -2322.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2615.8333
$ws.Range("I122").Value = 2339.2
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 7017.599999999999
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -4567.599999999999
$ws.Range("N122").Value = -16897
$ws.Range("H126").Value = 2174.8
$ws.Range("I126").Value = 2082.5557
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 6247.6671
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -3777.6671
$ws.Range("N126").Value = -13955
